# Add a "foaf:name" column (D) holding the full name, derived from the
# existing familyName/givenName columns in this Person sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell + new data cell.
$ws.Range("D1").Value = "foaf:name"
$ws.Range("D2").Value = "Minor Gordon"

# D1 (header) picks up the shared header font formatting.
$ws.Range("D1").Font.ThemeColor = 1

# D2 (data) picks up the formatting already used by the rest of row 2.
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
